# Refresh the Ccl3 -> Ccr5 LR-pair table with the new TPM-based NATMI run.
# The sending/target-cluster grid now spans all three clusters (ECs, MuSCs,
# Resolving-Mac) instead of the previous partial table, so rows 2-10 are
# rewritten in full (dimension grows from A1:T7 to A1:T10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row, columns A..T in sheet order:
#  Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
#  Ligand-expressing cells, Ligand detection rate, Ligand average/total
#  expression value, Ligand derived specificity (avg/total), Receptor-
#  expressing cells, Receptor detection rate, Receptor average/total
#  expression value, Receptor derived specificity (avg/total), Edge
#  average/total expression weight, Edge derived specificity (avg/total).
$rows = @(
    ,@("ECs", "Ccl3", "Ccr5", "ECs", 1, [double]"0.3333333333333333", [double]"0.07287233333333333", [double]"0.218617", [double]"0.0006119733110023554", [double]"0.0006119733110023554", 2, [double]"0.6666666666666666", [double]"0.022105", [double]"0.066315", [double]"0.0007557226718989593", [double]"0.0007557226718989592", [double]"0.001610842928333333", [double]"0.014497586355", [double]"4.624821057215529e-07", [double]"4.624821057215528e-07")
    ,@("ECs", "Ccl3", "Ccr5", "FAPs", 1, [double]"0.3333333333333333", [double]"0.07287233333333333", [double]"0.218617", [double]"0.0006119733110023554", [double]"0.0006119733110023554", 2, [double]"0.6666666666666666", [double]"0.4010506666666667", [double]"1.203152", [double]"0.01371106452749117", [double]"0.01371106452749117", [double]"0.02922549786488889", [double]"0.263029480784", [double]"8.390805556255715e-06", [double]"8.390805556255715e-06")
    ,@("ECs", "Ccl3", "Ccr5", "Resolving-Mac", 1, [double]"0.3333333333333333", [double]"0.07287233333333333", [double]"0.218617", [double]"0.0006119733110023554", [double]"0.0006119733110023554", 3, 1, [double]"28.82699233333333", [double]"86.480977", [double]"0.9855332128006099", [double]"0.9855332128006098", [double]"2.100690194312111", [double]"18.906211748809", [double]"0.0006031200233403782", [double]"0.0006031200233403781")
    ,@("MuSCs", "Ccl3", "Ccr5", "ECs", 1, [double]"0.3333333333333333", [double]"0.135447", [double]"0.406341", [double]"0.001137468024746511", [double]"0.001137468024746511", 2, [double]"0.6666666666666666", [double]"0.022105", [double]"0.066315", [double]"0.0007557226718989593", [double]"0.0007557226718989592", [double]"0.002994055935", [double]"0.026946503415", [double]"8.596103748610653e-07", [double]"8.596103748610651e-07")
    ,@("MuSCs", "Ccl3", "Ccr5", "FAPs", 1, [double]"0.3333333333333333", [double]"0.135447", [double]"0.406341", [double]"0.001137468024746511", [double]"0.001137468024746511", 2, [double]"0.6666666666666666", [double]"0.4010506666666667", [double]"1.203152", [double]"0.01371106452749117", [double]"0.01371106452749117", [double]"0.05432110964800001", [double]"0.488889986832", [double]"1.559589748525734e-05", [double]"1.559589748525734e-05")
    ,@("MuSCs", "Ccl3", "Ccr5", "Resolving-Mac", 1, [double]"0.3333333333333333", [double]"0.135447", [double]"0.406341", [double]"0.001137468024746511", [double]"0.001137468024746511", 3, 1, [double]"28.82699233333333", [double]"86.480977", [double]"0.9855332128006099", [double]"0.9855332128006098", [double]"3.904529630573", [double]"35.140766675157", [double]"0.001121012516886393", [double]"0.001121012516886393")
    ,@("Resolving-Mac", "Ccl3", "Ccr5", "ECs", 3, 1, [double]"118.8693136666667", [double]"356.607941", [double]"0.9982505586642512", [double]"0.9982505586642512", 2, [double]"0.6666666666666666", [double]"0.022105", [double]"0.066315", [double]"0.0007557226718989593", [double]"0.0007557226718989592", [double]"2.627606178601666", [double]"23.648455607415", [double]"0.0007544005794183768", [double]"0.0007544005794183767")
    ,@("Resolving-Mac", "Ccl3", "Ccr5", "FAPs", 3, 1, [double]"118.8693136666667", [double]"356.607941", [double]"0.9982505586642512", [double]"0.9982505586642512", 2, [double]"0.6666666666666666", [double]"0.4010506666666667", [double]"1.203152", [double]"0.01371106452749117", [double]"0.01371106452749117", [double]"47.67261749222578", [double]"429.053557430032", [double]"0.01368707782444965", [double]"0.01368707782444965")
    ,@("Resolving-Mac", "Ccl3", "Ccr5", "Resolving-Mac", 3, 1, [double]"118.8693136666667", [double]"356.607941", [double]"0.9982505586642512", [double]"0.9982505586642512", 3, 1, [double]"28.82699233333333", [double]"86.480977", [double]"0.9855332128006099", [double]"0.9855332128006098", [double]"3426.644793737595", [double]"30839.80314363835", [double]"0.9838090802603832", [double]"0.983809080260383")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $targetRow = $startRow + $i
    $rowValues = $rows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }
}
